# User checkpoint: Update portfolio presentation.
$p = $ppt.ActivePresentation

function Set-ParaText($textRange, $index, $newText) {
    # Clear first so the engine doesn't try to diff against the old
    # text and split the replacement across multiple runs.
    $para = $textRange.Paragraphs($index, 1)
    $para.Text = ""
    $para = $textRange.Paragraphs($index, 1)
    $para.Text = $newText
}

# --- Slide 1: update the date range in the subtitle ---
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
Set-ParaText $subtitle.TextFrame.TextRange 1 "2023.07 - 현재"

# --- Slide 3: rewrite the project detail / achievements bullets ---
$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

Set-ParaText $tr 1 "그로스폴리오 캠페인"
Set-ParaText $tr 3 "• OTT 서비스의 신규 구독자 유치를 위한 통합 마케팅 캠페인 기획 및 실행"
Set-ParaText $tr 4 "• 소셜미디어 채널별 맞춤형 콘텐츠 제작 및 광고 집행 전략 수립"
Set-ParaText $tr 5 "• 인플루언서 협업 프로그램 기획 및 운영을 통한 브랜드 인지도 확대"
Set-ParaText $tr 6 "• 주요 오리지널 콘텐츠 출시에 맞춘 시즌별 프로모션 캠페인 진행"
Set-ParaText $tr 7 "• 사용자 데이터 분석을 통한 타겟 맞춤형 리타겟팅 전략 수립"
Set-ParaText $tr 9 "• 캠페인 기간 중 신규 가입자 수 전월 대비 35% 증가"
Set-ParaText $tr 10 "• 소셜미디어 채널 팔로워 수 3개월간 25% 성장"
Set-ParaText $tr 11 "• 프로모션 참여율 평균 22% 달성 및 전환율 8.5% 기록"
